# Update "想去人数" (want-to-go count) values that changed between scrapes.
# Sheet "展览" (sheet1) rows: 2,5,6,8,9,11,12
# Sheet "全部类型" (sheet4) rows: 2,6,7,9,10,12,13

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 68
$wsExhibit.Range("F5").Value = 361
$wsExhibit.Range("F6").Value = 5465
$wsExhibit.Range("F8").Value = 5409
$wsExhibit.Range("F9").Value = 636
$wsExhibit.Range("F11").Value = 1390
$wsExhibit.Range("F12").Value = 25

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 68
$wsAll.Range("F6").Value = 361
$wsAll.Range("F7").Value = 5465
$wsAll.Range("F9").Value = 5409
$wsAll.Range("F10").Value = 636
$wsAll.Range("F12").Value = 1390
$wsAll.Range("F13").Value = 25
